$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H16").Value2 = 0
$ws.Range("I16").Value2 = 0
$ws.Range("K16").Value2 = 0
$ws.Range("M16").ClearContents()
$ws.Range("H137").Value2 = 2420331.8
$ws.Range("I137").Value2 = 1021346.1
$ws.Range("J137").Value2 = 7693431
$ws.Range("K137").Value2 = 3064038.3
$ws.Range("L137").Value2 = 23080293
$ws.Range("M137").Value2 = -3061488.3
$ws.Range("N137").Value2 = -23085393
$ws.Range("H141").Value2 = 1801.8628
$ws.Range("I141").Value2 = 1387.7941
$ws.Range("K141").Value2 = 4163.3823
$ws.Range("M141").Value2 = 1016.6177
# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H45").Value2 = 3120.963
$ws.Range("I45").Value2 = 2209.1052
$ws.Range("J45").Value2 = 5286.625
$ws.Range("K45").Value2 = 2209.1052
$ws.Range("L45").Value2 = 5286.625
$ws.Range("M45").Value2 = -1832.1052
$ws.Range("N45").Value2 = -6040.625
$ws.Range("H61").Value2 = 1436.3784
$ws.Range("I61").Value2 = 1487.7931
$ws.Range("J61").Value2 = 1250
$ws.Range("K61").Value2 = 1487.7931
$ws.Range("L61").Value2 = 1250
$ws.Range("M61").Value2 = -1275.7931
$ws.Range("N61").Value2 = -1674
$ws.Range("H132").Value2 = 106174.25
$ws.Range("I132").Value2 = 112408.086
$ws.Range("J132").Value2 = 12666.667
$ws.Range("K132").Value2 = 337224.258
$ws.Range("L132").Value2 = 38000.001
$ws.Range("M132").Value2 = -334694.258
$ws.Range("N132").Value2 = -43060.001
$ws.Range("H136").Value2 = 1436.3784
$ws.Range("I136").Value2 = 1487.7931
$ws.Range("J136").Value2 = 1250
$ws.Range("K136").Value2 = 4463.379300000001
$ws.Range("L136").Value2 = 3750
$ws.Range("M136").Value2 = -1913.379300000001
$ws.Range("N136").Value2 = -8850
# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value2 = 126241.04
$ws.Range("I134").Value2 = 159003
$ws.Range("J134").Value2 = 1745.6
$ws.Range("K134").Value2 = 477009
$ws.Range("L134").Value2 = 5236.799999999999
$ws.Range("M134").Value2 = -474474
$ws.Range("N134").Value2 = -10306.8
# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H52").Value2 = 28249.143
$ws.Range("I52").Value2 = 23800
$ws.Range("J52").Value2 = 28990.666
$ws.Range("K52").Value2 = 23800
$ws.Range("L52").Value2 = 28990.666
$ws.Range("M52").Value2 = -23506
$ws.Range("N52").Value2 = -29578.666
$ws.Range("H58").Value2 = 1864.9117
$ws.Range("I58").Value2 = 1988.5667
$ws.Range("J58").Value2 = 937.5
$ws.Range("K58").Value2 = 1988.5667
$ws.Range("L58").Value2 = 937.5
$ws.Range("M58").Value2 = -1785.5667
$ws.Range("N58").Value2 = -1343.5
$ws.Range("H92").Value2 = 29175
$ws.Range("J92").Value2 = 29175
$ws.Range("L92").Value2 = 29175
$ws.Range("N92").Value2 = -34167
$ws.Range("H132").Value2 = 3023.8718
$ws.Range("I132").Value2 = 2614.3794
$ws.Range("K132").Value2 = 7843.138199999999
$ws.Range("M132").Value2 = -5313.138199999999
$ws.Range("H134").Value2 = 6120.108
$ws.Range("I134").Value2 = 6758.4517
$ws.Range("J134").Value2 = 2822
$ws.Range("K134").Value2 = 20275.3551
$ws.Range("L134").Value2 = 8466
$ws.Range("M134").Value2 = -17740.3551
$ws.Range("N134").Value2 = -13536
$ws.Range("H136").Value2 = 1864.9117
$ws.Range("I136").Value2 = 1988.5667
$ws.Range("J136").Value2 = 937.5
$ws.Range("K136").Value2 = 5965.7001
$ws.Range("L136").Value2 = 2812.5
$ws.Range("M136").Value2 = -3415.7001
$ws.Range("N136").Value2 = -7912.5
# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H63").Value2 = 3919.5293
$ws.Range("I63").Value2 = 1999
$ws.Range("J63").Value2 = 5626.6665
$ws.Range("K63").Value2 = 5997
$ws.Range("L63").Value2 = 16879.9995
$ws.Range("M63").Value2 = -5248
$ws.Range("N63").Value2 = -18377.9995
$ws.Range("H64").Value2 = 5019.5757
$ws.Range("I64").Value2 = 823.7778
$ws.Range("J64").Value2 = 6593
$ws.Range("K64").Value2 = 2471.3334
$ws.Range("L64").Value2 = 19779
$ws.Range("M64").Value2 = -2201.3334
$ws.Range("N64").Value2 = -20319
$ws.Range("H66").Value2 = 3919.5293
$ws.Range("I66").Value2 = 1999
$ws.Range("J66").Value2 = 5626.6665
$ws.Range("K66").Value2 = 17991
$ws.Range("L66").Value2 = 50639.9985
$ws.Range("M66").Value2 = -14247
$ws.Range("N66").Value2 = -58127.9985
$ws.Range("H67").Value2 = 5019.5757
$ws.Range("I67").Value2 = 823.7778
$ws.Range("J67").Value2 = 6593
$ws.Range("K67").Value2 = 2471.3334
$ws.Range("L67").Value2 = 19779
$ws.Range("M67").Value2 = -1535.3334
$ws.Range("N67").Value2 = -21651
$ws.Range("H68").Value2 = 1950.3334
$ws.Range("I68").Value2 = 2378.5715
$ws.Range("J68").Value2 = 451.5
$ws.Range("K68").Value2 = 7135.7145
$ws.Range("L68").Value2 = 1354.5
$ws.Range("M68").Value2 = -6324.7145
$ws.Range("N68").Value2 = -2976.5
$ws.Range("H71").Value2 = 1950.3334
$ws.Range("I71").Value2 = 2378.5715
$ws.Range("J71").Value2 = 451.5
$ws.Range("K71").Value2 = 21407.1435
$ws.Range("L71").Value2 = 4063.5
$ws.Range("M71").Value2 = -17351.1435
$ws.Range("N71").Value2 = -12175.5
$ws.Range("H75").Value2 = 3778.6924
$ws.Range("I75").Value2 = 528.25
$ws.Range("J75").Value2 = 5223.3335
$ws.Range("K75").Value2 = 1584.75
$ws.Range("L75").Value2 = 15670.0005
$ws.Range("M75").Value2 = -586.75
$ws.Range("N75").Value2 = -17666.0005
$ws.Range("H76").Value2 = 7604.5454
$ws.Range("I76").Value2 = 3750
$ws.Range("J76").Value2 = 8461.111000000001
$ws.Range("K76").Value2 = 11250
$ws.Range("L76").Value2 = 25383.333
$ws.Range("M76").Value2 = -10867
$ws.Range("N76").Value2 = -26149.333
$ws.Range("H78").Value2 = 3778.6924
$ws.Range("I78").Value2 = 528.25
$ws.Range("J78").Value2 = 5223.3335
$ws.Range("K78").Value2 = 4754.25
$ws.Range("L78").Value2 = 47010.0015
$ws.Range("M78").Value2 = 237.75
$ws.Range("N78").Value2 = -56994.0015
$ws.Range("H79").Value2 = 7604.5454
$ws.Range("I79").Value2 = 3750
$ws.Range("J79").Value2 = 8461.111000000001
$ws.Range("K79").Value2 = 11250
$ws.Range("L79").Value2 = 25383.333
$ws.Range("M79").Value2 = -9924
$ws.Range("N79").Value2 = -28035.333
$ws.Range("H81").Value2 = 58825228
$ws.Range("I81").Value2 = 392.25
$ws.Range("J81").Value2 = 76925176
$ws.Range("K81").Value2 = 1176.75
$ws.Range("L81").Value2 = 230775528
$ws.Range("M81").Value2 = -53.75
$ws.Range("N81").Value2 = -230777774
$ws.Range("H82").Value2 = 7011.6665
$ws.Range("I82").Value2 = 752.4
$ws.Range("J82").Value2 = 8658.842000000001
$ws.Range("K82").Value2 = 2257.2
$ws.Range("L82").Value2 = 25976.526
$ws.Range("M82").Value2 = -1851.2
$ws.Range("N82").Value2 = -26788.526
$ws.Range("H84").Value2 = 58825228
$ws.Range("I84").Value2 = 392.25
$ws.Range("J84").Value2 = 76925176
$ws.Range("K84").Value2 = 3530.25
$ws.Range("L84").Value2 = 692326584
$ws.Range("M84").Value2 = 2085.75
$ws.Range("N84").Value2 = -692337816
$ws.Range("H85").Value2 = 7011.6665
$ws.Range("I85").Value2 = 752.4
$ws.Range("J85").Value2 = 8658.842000000001
$ws.Range("K85").Value2 = 2257.2
$ws.Range("L85").Value2 = 25976.526
$ws.Range("M85").Value2 = -853.1999999999998
$ws.Range("N85").Value2 = -28784.526
$ws.Range("H96").Value2 = 4562.5
$ws.Range("J96").Value2 = 4562.5
$ws.Range("L96").Value2 = 13687.5
$ws.Range("N96").Value2 = -17805.5
$ws.Range("H100").Value2 = 8054.5938
$ws.Range("J100").Value2 = 8054.5938
$ws.Range("L100").Value2 = 24163.7814
$ws.Range("N100").Value2 = -25785.7814
$ws.Range("H103").Value2 = 1292.5
$ws.Range("I103").Value2 = 350
$ws.Range("J103").Value2 = 1458.8235
$ws.Range("K103").Value2 = 1050
$ws.Range("L103").Value2 = 4376.470499999999
$ws.Range("M103").Value2 = -171
$ws.Range("N103").Value2 = -6134.470499999999
$ws.Range("H106").Value2 = 3749.5
$ws.Range("J106").Value2 = 3749.5
$ws.Range("L106").Value2 = 11248.5
$ws.Range("N106").Value2 = -13140.5
$ws.Range("H109").Value2 = 3372.7917
$ws.Range("I109").Value2 = 2246.2144
$ws.Range("J109").Value2 = 4950
$ws.Range("K109").Value2 = 6738.6432
$ws.Range("L109").Value2 = 14850
$ws.Range("M109").Value2 = -5698.6432
$ws.Range("N109").Value2 = -16930
$ws.Range("H112").Value2 = 20876938
$ws.Range("I112").Value2 = 1500
$ws.Range("J112").Value2 = 21784566
$ws.Range("K112").Value2 = 4500
$ws.Range("L112").Value2 = 65353698
$ws.Range("M112").Value2 = -3392
$ws.Range("N112").Value2 = -65355914
$ws.Range("H132").Value2 = 629863.8
$ws.Range("I132").Value2 = 1097511.6
$ws.Range("K132").Value2 = 9877604.4
$ws.Range("M132").Value2 = -9875074.4
$ws.Range("H137").Value2 = 1926.2963
$ws.Range("J137").Value2 = 3996
$ws.Range("L137").Value2 = 11988
$ws.Range("N137").Value2 = -22188
# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value2 = 0
$ws.Range("H132").Value2 = 1914.5366
$ws.Range("I132").Value2 = 1558.1538
$ws.Range("K132").Value2 = 4674.4614
$ws.Range("M132").Value2 = -2144.4614
# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H132").Value2 = 1597.2222
$ws.Range("I132").Value2 = 1517.4717
$ws.Range("K132").Value2 = 4552.4151
$ws.Range("M132").Value2 = -2022.4151
$ws.Range("H136").Value2 = 1444.9259
$ws.Range("I136").Value2 = 1344.0435
$ws.Range("J136").Value2 = 2025
$ws.Range("K136").Value2 = 4032.1305
$ws.Range("L136").Value2 = 6075
$ws.Range("M136").Value2 = -1482.1305
$ws.Range("N136").Value2 = -11175
# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H18").Value2 = 51005.25
$ws.Range("J18").Value2 = 51005.25
$ws.Range("L18").Value2 = 51005.25
$ws.Range("N18").Value2 = -51351.25
$ws.Range("H95").Value2 = 19268.8
$ws.Range("J95").Value2 = 19268.8
$ws.Range("L95").Value2 = 19268.8
$ws.Range("N95").Value2 = -24760.8
$ws.Range("H100").Value2 = 667452.75
$ws.Range("I100").Value2 = 526.8570999999999
$ws.Range("J100").Value2 = 1251012.9
$ws.Range("K100").Value2 = 1053.7142
$ws.Range("L100").Value2 = 2502025.8
$ws.Range("M100").Value2 = -512.7141999999999
$ws.Range("N100").Value2 = -2503107.8
$ws.Range("H132").Value2 = 2205.7856
$ws.Range("I132").Value2 = 2650.742
$ws.Range("J132").Value2 = 951.8182
$ws.Range("K132").Value2 = 7952.226000000001
$ws.Range("L132").Value2 = 2855.4546
$ws.Range("M132").Value2 = -5422.226000000001
$ws.Range("N132").Value2 = -7915.4546
$ws.Range("H136").Value2 = 1788.3556
$ws.Range("I136").Value2 = 1898.2778
$ws.Range("J136").Value2 = 1348.6666
$ws.Range("K136").Value2 = 5694.8334
$ws.Range("L136").Value2 = 4045.9998
$ws.Range("M136").Value2 = -3144.8334
$ws.Range("N136").Value2 = -9145.9998
$ws.Range("H137").Value2 = 27995
$ws.Range("J137").Value2 = 27995
$ws.Range("L137").Value2 = 27995
$ws.Range("N137").Value2 = -38195
